# codeforIATI SectorGroup.xlsx edit: the "category" and "group" columns were
# transposed. Column D (codeforiati:category-name) and column E
# (codeforiati:group-name) swap places, and column F (codeforiati:group-code)
# and column G (codeforiati:category-code) swap places, for every row
# (including the header row).
#
# The swap is done with Range.Cut (cut & paste) rather than copying
# .Value across, because plain value assignment lets Excel "smart" coerce
# numeric-looking text (e.g. "110") into a real number, which would change
# the cell type/style from the original text storage. Cut preserves the
# original cell content/type exactly, matching a real drag-and-drop style
# column swap in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$tmpRow = $lastRow + 1000

$dRange = $ws.Range("D1:D" + $lastRow)
$eRange = $ws.Range("E1:E" + $lastRow)
$tmpD = $ws.Range("D" + $tmpRow + ":D" + ($tmpRow + $lastRow - 1))

$dRange.Cut($tmpD)
$eRange.Cut($dRange)
$tmpD.Cut($eRange)

$fRange = $ws.Range("F1:F" + $lastRow)
$gRange = $ws.Range("G1:G" + $lastRow)
$tmpF = $ws.Range("F" + $tmpRow + ":F" + ($tmpRow + $lastRow - 1))

$fRange.Cut($tmpF)
$gRange.Cut($fRange)
$tmpF.Cut($gRange)
